$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new time-log entry row at 47 (shifts everything below down by 1)
$ws.Rows("47:47").Insert()

# Insert 3 more blank rows (49:51) to create spacing before the "Issues/Loose Ends" block
$ws.Rows("49:51").Insert()

# Fill in the new entry
$ws.Range("A47").Value = 43547
$ws.Range("B47").Value = 7.5
$ws.Range("D47").Value = "Week 7: got authentication working for admin and profile pages.  Tried to stop direct access via jsp but that didn't seem to work.`nFixed config so it doesn't kill my program.  Attempted a hibernate search - not getting expected result yet."

# Row heights to mirror the committed layout
$ws.Rows("47:47").RowHeight = 45
$ws.Rows("51:51").RowHeight = 14.25

# Update selection to match committed state
$ws.Range("D50").Select()
